# Add a "Registration" sheet after "Login", populate it by borrowing cell
# formatting from the existing Login sheet (so styles/borders match exactly
# without minting brand-new style records), then set the real values/hyperlink.

$wb = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item("Login")

# Put the selection on Login back to a full A1:E4 highlight (mirrors the
# state the workbook was saved in right after copying this range out).
$login.Activate()
$login.Range("A1:E4").Select()

# New sheet, placed immediately after "Login".
$reg = $wb.Worksheets.Add($null, $login)
$reg.Name = "Registration"

# Hyperlink style (s="4") must be applied before any cell copy lands on E2,
# otherwise adding the hyperlink later mints a duplicate style with
# applyFont="1" tacked on.
$reg.Hyperlinks.Add($reg.Range("E2"), "mailto:kiran@gmail.com") | Out-Null

# Row 1 (headers) - reuse Login's header formatting.
$login.Range("A1").Copy($reg.Range("A1"))   # tcId
$login.Range("B1").Copy($reg.Range("B1"))   # tcDescription (fill+border)
$login.Range("C1").Copy($reg.Range("C1"))   # -> name
$login.Range("D1").Copy($reg.Range("D1"))   # -> mobile
$login.Range("C1").Copy($reg.Range("E1"))   # -> eamil
$login.Range("D1").Copy($reg.Range("F1"))   # -> pwd
$login.Range("E1").Copy($reg.Range("G1"))   # expResult

$reg.Range("A1").Value = "tcId"
$reg.Range("B1").Value = "tcDescription"
$reg.Range("C1").Value = "name"
$reg.Range("D1").Value = "mobile"
$reg.Range("E1").Value = "eamil"
$reg.Range("F1").Value = "pwd"
$reg.Range("G1").Value = "expResult"

# Row 2 (single data row) - reuse Login's row4 (ValidInfo row) formatting.
$login.Range("A4").Copy($reg.Range("A2"))   # TCJBK03
$login.Range("B4").Copy($reg.Range("B2"))   # ValidInfo
$login.Range("C3").Copy($reg.Range("C2"))   # -> subhash (plain text style)
$login.Range("D4").Copy($reg.Range("D2"))   # -> mobile number (quotePrefix text style)
$login.Range("D4").Copy($reg.Range("F2"))   # -> pwd (quotePrefix text style)
$login.Range("E4").Copy($reg.Range("G2"))   # expResult text

$reg.Range("A2").Value = "TCJBK03"
$reg.Range("B2").Value = "ValidInfo"
$reg.Range("C2").Value = "subhash"
$reg.Range("D2").Value = "'9876545677"
$reg.Range("E2").Value = "kiran@gmail.com"
$reg.Range("F2").Value = "'123456"
$reg.Range("G2").Value = "User registered successfully."

Write-Output "done"
